# Ongoing testing of preemptive multitasking
#
# The "PC" saved-register row (old row 3) is removed from the virtualization
# worksheet; every row below it shifts up by one. After the shift, the
# formulas that used to reference the deleted row are repaired so the
# running totals in columns C and D keep working, and the D-column formula
# for the (new) PSP row is rewritten to reference C4 directly instead of
# continuing the old "previous-D-cell + current-B-cell" chain.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "PC" register row entirely - everything below shifts up one row.
$ws.Rows.Item(3).Delete()

# Row 3 is now the TOS row; its own C/D formulas were blown away by the
# deletion (they pointed at the row that no longer exists), so restore them.
$ws.Range("C3").Formula = "=B3-1"
$ws.Range("D3").Value = 0

# Rebuild the running totals in C4:C8 (NOS..ESP) as one fill so they come
# back as a proper shared formula, same as the original layout.
$ws.Range("C4:C8").Formula = "=C3+B4"

# Column D: D4 keeps the "previous D + current B" running total...
$ws.Range("D4").Formula = "=D3+B4"
# ...but D5 (PSP row) now instead reads off the C column directly.
$ws.Range("D5").Formula = "=C4+1"
# D6:D8 resume the normal "previous D + current B" running-total pattern.
$ws.Range("D6:D8").Formula = "=D5+B6"

# The selection left over from editing that formula is on D6.
$ws.Range("D6").Select()
